# Regenerate test data to include credit card account refund/credit transactions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Dining Out)
$ws.Range("C4").Value = -181022.72
$ws.Range("D4").Value = -181022.72

# Row 5 (Freelance Income)
$ws.Range("B5").Value = 338613.39
$ws.Range("D5").Value = 338613.39

# Row 6 (Groceries)
$ws.Range("C6").Value = -179006
$ws.Range("D6").Value = -179006

# Row 7 (Interest Income)
$ws.Range("B7").Value = 2254401.67
$ws.Range("D7").Value = 2254401.67

# Row 8 (Pets)
$ws.Range("C8").Value = -80719.32000000001
$ws.Range("D8").Value = -80719.32000000001

# Row 9 (Pharmacy)
$ws.Range("C9").Value = -190509.04
$ws.Range("D9").Value = -190509.04

# Row 10 (Rent)
$ws.Range("B10").Value = -80232.52
$ws.Range("D10").Value = -80232.52

# Row 11 (Shopping)
$ws.Range("C11").Value = -279746.85
$ws.Range("D11").Value = -279746.85

# Row 12 (Taxes)
$ws.Range("B12").Value = -141099.22
$ws.Range("D12").Value = -141099.22

# Row 13 (Transfer From)
$ws.Range("B13").Value = 54759.92
$ws.Range("C13").Value = 911003.9300000001
$ws.Range("D13").Value = 965763.85

# Row 14 (Transfer To)
$ws.Range("B14").Value = -965763.85
$ws.Range("D14").Value = -965763.85

# Row 15 (Utilities)
$ws.Range("B15").Value = -83002.44
$ws.Range("D15").Value = -83002.44

# Row 16 (Wages & Salary)
$ws.Range("B16").Value = 413324.54
$ws.Range("D16").Value = 413324.54

# Row 17 (Total)
$ws.Range("B17").Value = 1791001.49
$ws.Range("D17").Value = 1791001.49
